$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.124600112295553
$ws.Range("C2").Value = 0.03465183341573663
$ws.Range("D2").Value = 0.003026854328761708
$ws.Range("E2").Value = 0.06674794700219699
$ws.Range("F2").Value = 4.687133108435773
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1815652751018284
$ws.Range("K2").Value = 1.585180651292916
$ws.Range("L2").Value = 0.2706623424018915
$ws.Range("M2").Value = 0.4434864055175858
$ws.Range("N2").Value = 4.954385763929537
$ws.Range("B3").Value = 2.095091171269956
$ws.Range("C3").Value = 0.03053409699747078
$ws.Range("D3").Value = 0.003034667849737716
$ws.Range("E3").Value = 0.06694640041787991
$ws.Range("F3").Value = 4.678831105747776
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1820829666192836
$ws.Range("K3").Value = 1.553133237958463
$ws.Range("L3").Value = 0.2701311971990705
$ws.Range("M3").Value = 0.4391798820445487
$ws.Range("N3").Value = 4.962651873587646
$ws.Range("B4").Value = 2.078074266701861
$ws.Range("C4").Value = 0.02801690480386299
$ws.Range("D4").Value = 0.00304169749619243
$ws.Range("E4").Value = 0.06708069524931926
$ws.Range("F4").Value = 4.675267130709813
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1824263713628049
$ws.Range("K4").Value = 1.534345845640274
$ws.Range("L4").Value = 0.2699131463232547
$ws.Range("M4").Value = 0.4367519697315636
$ws.Range("N4").Value = 4.968628249207143
$ws.Range("B5").Value = 2.071417062371722
$ws.Range("C5").Value = 0.02699386853315389
$ws.Range("D5").Value = 0.003045128267902264
$ws.Range("E5").Value = 0.06713855970862159
$ws.Range("F5").Value = 4.674200608056225
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1825727468248264
$ws.Range("K5").Value = 1.526913707004383
$ws.Range("L5").Value = 0.2698515205998575
$ws.Range("M5").Value = 0.4358170429632011
$ws.Range("N5").Value = 4.971290242682272
$ws.Range("B6").Value = 2.07032839697186
$ws.Range("C6").Value = 0.02682415803369054
$ws.Range("D6").Value = 0.003045732282761726
$ws.Range("E6").Value = 0.06714835786951934
$ws.Range("F6").Value = 4.67404682231448
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1825974414142628
$ws.Range("K6").Value = 1.52569313123567
$ws.Range("L6").Value = 0.2698429340734805
$ws.Range("M6").Value = 0.4356650909848305
$ws.Range("N6").Value = 4.971745949734725
$ws.Range("B7").Value = 2.077983362062525
$ws.Range("C7").Value = 0.02800309675703261
$ws.Range("D7").Value = 0.003041741466232928
$ws.Range("E7").Value = 0.06708146291069461
$ws.Range("F7").Value = 4.675251184742024
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1824283193602394
$ws.Range("K7").Value = 1.534244706596382
$ws.Range("L7").Value = 0.269912204882587
$ws.Range("M7").Value = 0.436739140347747
$ws.Range("N7").Value = 4.96866323239631
$ws.Range("B8").Value = 2.114196927825901
$ws.Range("C8").Value = 0.03322969673264708
$ws.Range("D8").Value = 0.003029088013441594
$ws.Range("E8").Value = 0.066813796540802
$ws.Range("F8").Value = 4.683952420569895
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1817384820798758
$ws.Range("K8").Value = 1.573946068750018
$ws.Range("L8").Value = 0.2704568051335769
$ws.Range("M8").Value = 0.4419566626578018
$ws.Range("N8").Value = 4.957048935011912
$ws.Range("B9").Value = 2.193947096612362
$ws.Range("C9").Value = 0.04357022541930178
$ws.Range("D9").Value = 0.003021782675126872
$ws.Range("E9").Value = 0.06638724160044251
$ws.Range("F9").Value = 4.713178617154426
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1805877909151974
$ws.Range("K9").Value = 1.658862247266512
$ws.Range("L9").Value = 0.2723801896238527
$ws.Range("M9").Value = 0.4539023922264889
$ws.Range("N9").Value = 4.941422483043809
$ws.Range("B10").Value = 2.257866694958523
$ws.Range("C10").Value = 0.0512278868410192
$ws.Range("D10").Value = 0.003026831374881134
$ws.Range("E10").Value = 0.0661332772248846
$ws.Range("F10").Value = 4.742067719299726
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1798648009064632
$ws.Range("K10").Value = 1.725566787454426
$ws.Range("L10").Value = 0.2743124820923484
$ws.Range("M10").Value = 0.463722706681736
$ws.Range("N10").Value = 4.934302829133259
$ws.Range("B11").Value = 2.288103256788702
$ws.Range("C11").Value = 0.05472578712610243
$ws.Range("D11").Value = 0.003031338556683139
$ws.Range("E11").Value = 0.06603053666038505
$ws.Range("F11").Value = 4.75682157244708
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1795623159598598
$ws.Range("K11").Value = 1.756852895071546
$ws.Range("L11").Value = 0.2753038214119385
$ws.Range("M11").Value = 0.4684167409202331
$ws.Range("N11").Value = 4.932011715581297
$ws.Range("B12").Value = 2.299719653436853
$ws.Range("C12").Value = 0.05605249330537276
$ws.Range("D12").Value = 0.00303335903619395
$ws.Range("E12").Value = 0.06599346164794451
$ws.Range("F12").Value = 4.762640239423661
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1794515574381457
$ws.Range("K12").Value = 1.768835643942452
$ws.Range("L12").Value = 0.2756953232585815
$ws.Range("M12").Value = 0.4702268103815115
$ws.Range("N12").Value = 4.93128045146716
$ws.Range("B13").Value = 2.297210454450521
$ws.Range("C13").Value = 0.05576666785709961
$ws.Range("D13").Value = 0.003032910002437106
$ws.Range("E13").Value = 0.06600136513220978
$ws.Range("F13").Value = 4.761376781971052
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1794752430650881
$ws.Range("K13").Value = 1.766248924478106
$ws.Range("L13").Value = 0.2756102909237441
$ws.Range("M13").Value = 0.4698355332598325
$ws.Range("N13").Value = 4.931431877923501
$ws.Range("B14").Value = 2.289055609761874
$ws.Range("C14").Value = 0.05483489307628986
$ws.Range("D14").Value = 0.00303149851923834
$ws.Range("E14").Value = 0.0660274498367901
$ws.Range("F14").Value = 4.75729563460402
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1795531279776093
$ws.Range("K14").Value = 1.757836010170337
$ws.Range("L14").Value = 0.2753357080559695
$ws.Range("M14").Value = 0.4685650046134384
$ws.Range("N14").Value = 4.931948821413158
$ws.Range("B15").Value = 2.28408220421602
$ws.Range("C15").Value = 0.05426443287086613
$ws.Range("D15").Value = 0.00303067467873408
$ws.Range("E15").Value = 0.06604366562147757
$ws.Range("F15").Value = 4.754825985474525
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1796013274462176
$ws.Range("K15").Value = 1.752700485839966
$ws.Range("L15").Value = 0.275169613762742
$ws.Range("M15").Value = 0.4677910050555667
$ws.Range("N15").Value = 4.932283220370778
$ws.Range("B16").Value = 2.255913971357757
$ws.Range("C16").Value = 0.05099958535633675
$ws.Range("D16").Value = 0.003026580941432933
$ws.Range("E16").Value = 0.06614024818824848
$ws.Range("F16").Value = 4.741135949515154
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1798850994131822
$ws.Range("K16").Value = 1.723541114289617
$ws.Range("L16").Value = 0.2742499518342285
$ws.Range("M16").Value = 0.4634204974346048
$ws.Range("N16").Value = 4.934471632974208
$ws.Range("B17").Value = 2.238930412143645
$ws.Range("C17").Value = 0.04900043417235622
$ws.Range("D17").Value = 0.003024632721278842
$ws.Range("E17").Value = 0.06620276802613567
$ws.Range("F17").Value = 4.733150362289336
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1800659396665125
$ws.Range("K17").Value = 1.705893975698473
$ws.Range("L17").Value = 0.2737145038056497
$ws.Range("M17").Value = 0.460797357428973
$ws.Range("N17").Value = 4.936056910807224
$ws.Range("B18").Value = 2.229271042079461
$ws.Range("C18").Value = 0.04785192642617631
$ws.Range("D18").Value = 0.003023720459800927
$ws.Range("E18").Value = 0.0662399320552316
$ws.Range("F18").Value = 4.728709003989522
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.180172440401849
$ws.Range("K18").Value = 1.695832465619532
$ws.Range("L18").Value = 0.2734171042082352
$ws.Range("M18").Value = 0.4593099406316696
$ws.Range("N18").Value = 4.93705791432167
$ws.Range("B19").Value = 2.226019294536115
$ws.Range("C19").Value = 0.04746329174365371
$ws.Range("D19").Value = 0.003023447488236464
$ws.Range("E19").Value = 0.06625272225026269
$ws.Range("F19").Value = 4.727231301377316
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1802089270757516
$ws.Range("K19").Value = 1.692441040157263
$ws.Range("L19").Value = 0.2733182277684705
$ws.Range("M19").Value = 0.458809994619827
$ws.Range("N19").Value = 4.937412154290755
$ws.Range("B20").Value = 2.24072704875482
$ws.Range("C20").Value = 0.04921310697359615
$ws.Range("D20").Value = 0.003024818584154154
$ws.Range("E20").Value = 0.06619598810656768
$ws.Range("F20").Value = 4.733984737928139
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1800464316817418
$ws.Range("K20").Value = 1.707763368498576
$ws.Range("L20").Value = 0.2737704089546114
$ws.Range("M20").Value = 0.4610743865949303
$ws.Range("N20").Value = 4.935878923628636
$ws.Range("B21").Value = 2.291446369232062
$ws.Range("C21").Value = 0.05510851972036335
$ws.Range("D21").Value = 0.003031904624647197
$ws.Range("E21").Value = 0.06601973851139586
$ws.Range("F21").Value = 4.758488078987995
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1795301486302616
$ws.Range("K21").Value = 1.76030341335138
$ws.Range("L21").Value = 0.2754159230585174
$ws.Range("M21").Value = 0.4689373070211005
$ws.Range("N21").Value = 4.931793282100301
$ws.Range("B22").Value = 2.32556457349483
$ws.Range("C22").Value = 0.05897393911813253
$ws.Range("D22").Value = 0.003038362460181965
$ws.Range("E22").Value = 0.06591521553868684
$ws.Range("F22").Value = 4.775852926325769
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1792147901325443
$ws.Range("K22").Value = 1.79543043764852
$ws.Range("L22").Value = 0.2765852005556013
$ws.Range("M22").Value = 0.4742658235654744
$ws.Range("N22").Value = 4.929917750657665
$ws.Range("B23").Value = 2.307266345851986
$ws.Range("C23").Value = 0.0569097358367685
$ws.Range("D23").Value = 0.003034749983902962
$ws.Range("E23").Value = 0.06597002821846676
$ws.Range("F23").Value = 4.766461443798335
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1793810879014828
$ws.Range("K23").Value = 1.776610302738618
$ws.Range("L23").Value = 0.2759525651531121
$ws.Range("M23").Value = 0.4714045617621707
$ws.Range("N23").Value = 4.930846023776837
$ws.Range("B24").Value = 2.239914463175751
$ws.Range("C24").Value = 0.04911695501776592
$ws.Range("D24").Value = 0.003024733908808308
$ws.Range("E24").Value = 0.06619904950608468
$ws.Range("F24").Value = 4.733607050475101
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1800552433498854
$ws.Range("K24").Value = 1.706917954151891
$ws.Range("L24").Value = 0.2737451017244581
$ws.Range("M24").Value = 0.4609490773455747
$ws.Range("N24").Value = 4.935959112476198
$ws.Range("B25").Value = 2.171437214005067
$ws.Range("C25").Value = 0.04076257293232288
$ws.Range("D25").Value = 0.003021911779153008
$ws.Range("E25").Value = 0.06649216269002967
$ws.Range("F25").Value = 4.703969869793013
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1808775294856844
$ws.Range("K25").Value = 1.635132979706384
$ws.Range("L25").Value = 0.2717685249063564
$ws.Range("M25").Value = 0.4504873723806995
$ws.Range("N25").Value = 4.944884215423215
